# Sheet1 originally had a schedule table (A:№, B:Посыл/Заповедь, C:Время)
# spanning rows 1-11. Rows 7-11 all shared the same "Посыл" text (the long
# encrypted blob) with consecutive 5-minute time slots (22:30-22:35 ...
# 22:50-22:55). The edit collapses those five rows into a single row 7:
#   B7 -> the short "Посыл" text (same one already used in rows 3 & 5)
#   C7 -> "19:0-19:5"
# and removes the now-redundant rows 8-11, shrinking the used range from
# A1:C11 down to A1:C7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the replacement "Посыл" text from a row that already has it before
# we start deleting anything.
$shortPosyl = $ws.Range("B3").Value()

# Drop the four trailing rows (old rows 8-11); row 7 survives and gets
# overwritten below with the new combined content.
$ws.Rows("8:11").Delete() | Out-Null

$ws.Range("B7").Value = $shortPosyl
$ws.Range("C7").Value = "19:0-19:5"

# Match the author's final selection (was B15, now B13 after the rows moved up).
$ws.Range("B13").Select() | Out-Null
